# Automatische test-sync: 2025-07-29 21:57:50
# Adds Testmail #14 ("Heb je de CE-certificaten van dit product?") as row 16
# of the "Logs" sheet, extends the conditional-formatting ranges to cover it,
# and updates the "Dashboard" summary sheet so "Productinformatie" (now 5
# mails) sorts above "Overig" (still 4 mails).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet - append the new test-mail row
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(16, 1).Value = "Heb je de CE-certificaten van dit product?"
$logs.Cells.Item(16, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(16, 3).Value = "Testmail #14: Heb je de CE-certificaten van dit product?"
$logs.Cells.Item(16, 4).Value = "Productinformatie"
$logs.Cells.Item(16, 5).Value = "Beste klant," + "`n" + "Dank u voor uw e-mail. Wij kunnen u bevestigen dat dit product over de vereiste CE-certificaten beschikt. Mocht u nog verdere vragen hebben of meer informatie nodig hebben, aarzel dan niet om contact met ons op te nemen." + "`n" + "Met vriendelijke groet," + "`n" + "[Naam bedrijf] E-mailassistent"
$logs.Cells.Item(16, 6).Value = "2025-07-29 21:57:40"
$logs.Cells.Item(16, 7).Value = "Ja"
$logs.Cells.Item(16, 8).Value = "Nee"
$logs.Cells.Item(16, 9).Value = "Ja"
$logs.Cells.Item(16, 10).Value = "Nee"

# The multi-line text in E16 makes the engine auto-expand the row height;
# AutoFit()-ing the row back collapses it to the (implicit) default height
# again, the same as every other data row on this sheet.
$logs.Rows.Item(16).EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 2. Grow every conditional-formatting rule on "Logs" from row 15 to row 16
# ---------------------------------------------------------------------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range($col + "2:" + $col + "15")
    $newRange = $logs.Range($col + "2:" + $col + "16")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count(); $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. "Dashboard" sheet - Productinformatie now leads with 5, Overig stays 4
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(2, 1).Value = "Productinformatie"
$dash.Cells.Item(2, 2).Value = 5
$dash.Cells.Item(3, 1).Value = "Overig"
$dash.Cells.Item(3, 2).Value = 4
